$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.112.42"
$ws.Range("E2").Value = "  -4.07%  "
$ws.Range("D3").Value = "2.250.55"
$ws.Range("E3").Value = "  -4.46%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "233.47"
$ws.Range("E5").Value = "  -3.11%  "
$ws.Range("D6").Value = "0.635"
$ws.Range("E6").Value = "  -6.13%  "
$ws.Range("D7").Value = "70.15"
$ws.Range("E7").Value = "  -4.59%  "
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("D9").Value = "0.558"
$ws.Range("E9").Value = "  -7.24%  "
$ws.Range("D10").Value = "0.0992"
$ws.Range("E10").Value = "  -1.51%  "
$ws.Range("D11").Value = "58.41"
$ws.Range("E11").Value = "  -1.24%  "
$ws.Range("D12").Value = "36.09"
$ws.Range("E12").Value = "  +7.84%  "
$ws.Range("E13").Value = "  -3.45%  "
$ws.Range("D14").Value = "6.81"
$ws.Range("E14").Value = "  -6.96%  "
$ws.Range("D15").Value = "2.586.69"
$ws.Range("E15").Value = "  -4.45%  "
$ws.Range("D16").Value = "15.09"
$ws.Range("E16").Value = "  -8.35%  "
$ws.Range("D17").Value = "0.865"
$ws.Range("E17").Value = "  -4.88%  "
$ws.Range("D18").Value = "2.256.70"
$ws.Range("E18").Value = "  -4.26%  "
$ws.Range("D19").Value = "42.041.26"
$ws.Range("E19").Value = "  -4.06%  "
$ws.Range("D20").Value = "0.0₃0979"
$ws.Range("E20").Value = "  -4.68%  "
$ws.Range("D21").Value = "6.27"
$ws.Range("E21").Value = "  -6.51%  "
$ws.Range("D22").Value = "73.50"
$ws.Range("E22").Value = "  -5.28%  "
$ws.Range("D23").Value = "237.71"
$ws.Range("E23").Value = "  -7.50%  "
$ws.Range("D24").Value = "2.07"
$ws.Range("E24").Value = "  +5.96%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").Value = "3.65"
$ws.Range("E26").Value = "  -3.24%  "
$ws.Range("D27").Value = "2.36"
$ws.Range("E27").Value = "  -5.62%  "
$ws.Range("D28").Value = "10.06"
$ws.Range("E28").Value = "  -5.28%  "
$ws.Range("D29").Value = "2.18"
$ws.Range("E29").Value = "  -0.24%  "
$ws.Range("D30").Value = "169.50"
$ws.Range("E30").Value = "  -4.64%  "
$ws.Range("D31").Value = "20.70"
$ws.Range("E31").Value = "  -8.97%  "
$ws.Range("E32").Value = "  -6.97%  "
$ws.Range("D33").Value = "0.128"
$ws.Range("E33").Value = "  -6.77%  "
$ws.Range("D34").Value = "5.49"
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("D35").Value = "0.0721"
$ws.Range("E35").Value = "  -4.76%  "
$ws.Range("E36").Value = "  -8.01%  "
$ws.Range("D37").Value = "3.62"
$ws.Range("E37").Value = "  -5.02%  "
$ws.Range("D38").Value = "22.13"
$ws.Range("E38").Value = "  +16.51%  "
$ws.Range("E39").Value = "  -4.93%  "
$ws.Range("D40").Value = "6.07"
$ws.Range("E40").Value = "  -5.98%  "
$ws.Range("E41").Value = "  -3.73%  "
$ws.Range("D42").Value = "67.33"
$ws.Range("E42").Value = "  -1.30%  "
$ws.Range("B43").Value = "FTXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D43").Value = "5.03"
$ws.Range("E43").Value = "  -3.00%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "9.15"
$ws.Range("E44").Value = "  -2.25%  "
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").Value = "0.193"
$ws.Range("E45").Value = "  -4.98%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").Value = "0.101"
$ws.Range("E46").Value = "  -8.80%  "
$ws.Range("E47").Value = "  +0.21%  "
$ws.Range("E48").Value = "  -5.47%  "
$ws.Range("D49").Value = "4.36"
$ws.Range("E49").Value = "  +6.95%  "
$ws.Range("D50").Value = "1.18"
$ws.Range("E50").Value = "  -6.47%  "
$ws.Range("D51").Value = "9.92"
$ws.Range("E51").Value = "  +2.56%  "
